# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Mon May  6 01:54:50 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.358.69'
$ws.Range("E2").Value = '  +1.72%  '

# Row 3
$ws.Range("D3").Value = '3.157.24'
$ws.Range("E3").Value = '  +2.05%  '

# Row 4
$ws.Range("E4").Value = '  +0.34%  '

# Row 5
$ws.Range("D5").Value = '''592.60'
$ws.Range("E5").Value = '  +1.66%  '

# Row 6
$ws.Range("D6").Value = '''146.85'
$ws.Range("E6").Value = '  +1.67%  '

# Row 7
$ws.Range("E7").Value = '  +0.12%  '

# Row 8
$ws.Range("D8").Value = '3.144.61'
$ws.Range("E8").Value = '  +1.86%  '

# Row 9
$ws.Range("E9").Value = '  +0.86%  '

# Row 10
$ws.Range("E10").Value = '  +4.15%  '

# Row 11
$ws.Range("D11").Value = '''5.90'
$ws.Range("E11").Value = '  +4.00%  '

# Row 12
$ws.Range("D12").Value = '''0.458'
$ws.Range("E12").Value = '  +0.70%  '

# Row 13
$ws.Range("D13").Value = '''0.0000250'
$ws.Range("E13").Value = '  +1.55%  '

# Row 14
$ws.Range("D14").Value = '''37.32'
$ws.Range("E14").Value = '  -0.18%  '

# Row 15
$ws.Range("D15").Value = '3.679.49'
$ws.Range("E15").Value = '  +2.04%  '

# Row 16
$ws.Range("E16").Value = '  -0.32%  '

# Row 17
$ws.Range("D17").Value = '''7.28'
$ws.Range("E17").Value = '  +2.36%  '

# Row 18
$ws.Range("D18").Value = '64.135.26'
$ws.Range("E18").Value = '  +1.52%  '

# Row 19
$ws.Range("D19").Value = '3.149.66'
$ws.Range("E19").Value = '  +1.94%  '

# Row 20
$ws.Range("D20").Value = '''468.45'
$ws.Range("E20").Value = '  +1.98%  '

# Row 21
$ws.Range("D21").Value = '''14.39'
$ws.Range("E21").Value = '  +0.72%  '

# Row 22
$ws.Range("D22").Value = '''0.735'
$ws.Range("E22").Value = '  +1.45%  '

# Row 23
$ws.Range("D23").Value = '''7.60'
$ws.Range("E23").Value = '  +2.48%  '

# Row 24
$ws.Range("D24").Value = '''2.40'
$ws.Range("E24").Value = '  +12.99%  '

# Row 25
$ws.Range("D25").Value = '''13.24'
$ws.Range("E25").Value = '  +1.79%  '

# Row 26
$ws.Range("D26").Value = '''81.42'
$ws.Range("E26").Value = '  +0.87%  '

# Row 27
$ws.Range("E27").Value = '  +0.05%  '

# Row 28
$ws.Range("D28").Value = '''9.86'
$ws.Range("E28").Value = '  +11.47%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.73'
$ws.Range("E29").Value = '  +2.42%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '''2.24'
$ws.Range("E30").Value = '  +2.15%  '

# Row 31
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''7.36'
$ws.Range("E31").Value = '  +8.34%  '

# Row 32
$ws.Range("E32").Value = '  +0.29%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.115'
$ws.Range("E33").Value = '  +8.03%  '

# Row 34
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '''27.71'
$ws.Range("E34").Value = '  +3.96%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0867'
$ws.Range("E35").Value = '  +3.06%  '

# Row 36
$ws.Range("E36").Value = '  +3.66%  '

# Row 37
$ws.Range("D37").Value = '''6.17'
$ws.Range("E37").Value = '  +2.84%  '

# Row 38
$ws.Range("D38").Value = '''2.30'
$ws.Range("E38").Value = '  +0.08%  '

# Row 39
$ws.Range("D39").Value = '''3.23'
$ws.Range("E39").Value = '  -2.57%  '

# Row 40
$ws.Range("D40").Value = '''467.30'
$ws.Range("E40").Value = '  +7.50%  '

# Row 41
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = '''9.44'
$ws.Range("E41").Value = '  +8.68%  '

# Row 42
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''51.31'
$ws.Range("E42").Value = '  +1.93%  '

# Row 43
$ws.Range("D43").Value = '''0.294'
$ws.Range("E43").Value = '  +8.51%  '

# Row 44
$ws.Range("D44").Value = '''0.0375'
$ws.Range("E44").Value = '  +2.11%  '

# Row 45
$ws.Range("D45").Value = '2.918.10'
$ws.Range("E45").Value = '  +1.93%  '

# Row 46
$ws.Range("D46").Value = '''40.19'
$ws.Range("E46").Value = '  +13.38%  '

# Row 47
$ws.Range("E47").Value = '  -0.27%  '

# Row 48
$ws.Range("D48").Value = '''133.99'
$ws.Range("E48").Value = '  +8.25%  '

# Row 49
$ws.Range("E49").Value = '  +0.01%  '

# Row 50
$ws.Range("D50").Value = '''2.25'
$ws.Range("E50").Value = '  +5.22%  '

# Row 51
$ws.Range("E51").Value = '  +1.15%  '
